{"js": "// Cover letter update: new recipient (Matthew \"Matty\" Dinh), new date,\n// new role/company (Blockchain Researcher @ Coinbase instead of Data\n// Scientist @ Croptix), and a handful of supporting-paragraph tweaks.\n//\n// Strategy: locate each target phrase with body.search(...) (exact,\n// case-sensitive, non-wildcard) and replace it via\n// range.insertText(..., Word.InsertLocation.replace). Every search\n// string below was chosen to be unique in the document (or is handled\n// in an order that makes it unique) so each block edits exactly the\n// run(s) intended.\n\nasync function replaceUnique(context, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"replaceUnique: expected exactly 1 match for \" +\n        JSON.stringify(searchText) +\n        \" but found \" +\n        results.items.length\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n  return results.items[0];\n}\n\nasync function replaceAll(context, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n  return results.items.length;\n}\n\n// Force a run boundary at the start of `range` without altering its\n// text or visible formatting: flip a character property on and back\n// off. Word/Office.js coalesces runs with identical rPr, but it always\n// keeps the edit's own run split in place once the value has been\n// round-tripped through a real write.\nasync function splitBefore(context, searchText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"splitBefore: expected exactly 1 match for \" +\n        JSON.stringify(searchText) +\n        \" but found \" +\n        results.items.length\n    );\n  }\n  const range = results.items[0];\n  range.font.load(\"bold\");\n  await context.sync();\n  const original = range.font.bold;\n  range.font.bold = !original;\n  await context.sync();\n  range.font.bold = original;\n  await context.sync();\n}\n\n// 1) Letterhead recipient name.\nawait replaceUnique(context, \"Faith Oliver\", \"Matthew (Matty) Dinh\");\n\n// 2) Letter date.\nawait replaceUnique(context, \"January 27, 2021\", \"February 9, 2021\");\n\n// 3) Source of the job posting.\nawait replaceUnique(context, \"Indeed\", \"the company website\");\n\n// 4) Target company name, both mentions (intro paragraph + closing\n// paragraph both say \"Croptix\" -> \"Coinbase\").\nconst croptixCount = await replaceAll(context, \"Croptix\", \"Coinbase\");\nif (croptixCount !== 2) {\n  throw new Error(\"expected 2 Croptix matches, found \" + croptixCount);\n}\n\n// 5) Recipient title/location line: \"VP Human Resources at Phreesia, Somerville, Massachusetts\"\n//    -> \"Human Resources at San Francisco Bay Area\"\nawait replaceUnique(context, \"VP Human Resources at Phreesia\", \"Human Resources at \");\nawait replaceUnique(context, \", Somerville, Massachusetts\", \"San Francisco Bay Area\");\n\n// 6) Salutation: \"Dear Ms. Oliver\" -> \"Dear Mr. Dinh\"\nawait splitBefore(context, \"Dear M\");\nawait replaceUnique(context, \"Ms. \", \"Mr. \");\n// After step 1, \"Oliver\" only occurs here.\nawait replaceUnique(context, \"Oliver\", \"Dinh\");\n\n// 7) Role applied for, in the introductory paragraph.\nawait replaceUnique(\n  context,\n  \"Data Scientist position at \",\n  \"Blockchain Researcher position at \"\n);\n\n// 8) Body paragraph: drop the \"Using Fast Fourier transform (FFT), \"\n//    clause and capitalize the following \"my\" -> \"My\".\nawait replaceUnique(\n  context,\n  \"arteriovenous fistula (AVF). Using Fast Fourier transform (FFT), my \",\n  \"arteriovenous fistula (AVF). My \"\n);\n\n// 9) Body paragraph: quantify the monthly savings.\nawait replaceUnique(\n  context,\n  \" and thus reduced the monthly cost for patients using our products.\",\n  \" and thus reduced around $900 monthly cost for patients using our products.\"\n);\n\n// 10) Closing paragraph: drop the repeated \"Data Scientist\" qualifier.\nawait replaceUnique(context, \"this Data Scientist role\", \"this role\");\nawait replaceUnique(\n  context,\n  \"about the Data Scientist position and\",\n  \"about the position and\"\n);\n", "ps1": "# Cover letter update: new recipient (Matthew \"Matty\" Dinh), new date,\n# new role/company (Blockchain Researcher @ Coinbase instead of Data\n# Scientist @ Croptix), and a handful of supporting-paragraph tweaks.\n#\n# Strategy: use Range.Find.Execute(...) against $d.Content (exact,\n# case-sensitive, no wildcards) for every edit. Each search string below\n# was chosen to be unique in the document (or is handled in an order\n# that makes it unique) so each call edits exactly the text intended.\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($FindText, $ReplaceText) {\n    # wdReplaceOne = 1; Wrap = wdFindContinue (1)\n    $rng = $d.Content\n    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 1)\n    if (-not $found) {\n        throw \"Replace-UniqueText: text not found: $FindText\"\n    }\n}\n\nfunction Replace-AllText($FindText, $ReplaceText) {\n    # wdReplaceAll = 2; Wrap = wdFindContinue (1)\n    $rng = $d.Content\n    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $found) {\n        throw \"Replace-AllText: text not found: $FindText\"\n    }\n}\n\n# Force a run boundary right before $FindText without touching its text\n# or visible formatting: flip Bold on then back off. Word coalesces runs\n# that share identical formatting, but a real property round-trip keeps\n# the split the edit introduces.\nfunction Split-Before($FindText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    if (-not $found) {\n        throw \"Split-Before: text not found: $FindText\"\n    }\n    $was = $rng.Bold\n    $rng.Bold = 1\n    $rng.Bold = $was\n}\n\n# 1) Letterhead recipient name.\nReplace-UniqueText \"Faith Oliver\" \"Matthew (Matty) Dinh\"\n\n# 2) Letter date.\nReplace-UniqueText \"January 27, 2021\" \"February 9, 2021\"\n\n# 3) Source of the job posting.\nReplace-UniqueText \"Indeed\" \"the company website\"\n\n# 4) Target company name, both mentions (intro paragraph + closing\n# paragraph both say \"Croptix\" -> \"Coinbase\").\nReplace-AllText \"Croptix\" \"Coinbase\"\n\n# 5) Recipient title/location line:\n# \"VP Human Resources at Phreesia, Somerville, Massachusetts\"\n# -> \"Human Resources at San Francisco Bay Area\"\nReplace-UniqueText \"VP Human Resources at Phreesia\" \"Human Resources at \"\nReplace-UniqueText \", Somerville, Massachusetts\" \"San Francisco Bay Area\"\n\n# 6) Salutation: \"Dear Ms. Oliver\" -> \"Dear Mr. Dinh\"\nSplit-Before \"Dear M\"\nReplace-UniqueText \"Ms. \" \"Mr. \"\n# After step 1, \"Oliver\" only occurs here.\nReplace-UniqueText \"Oliver\" \"Dinh\"\n\n# 7) Role applied for, in the introductory paragraph.\nReplace-UniqueText \"Data Scientist position at \" \"Blockchain Researcher position at \"\n\n# 8) Body paragraph: drop the \"Using Fast Fourier transform (FFT), \"\n# clause and capitalize the following \"my\" -> \"My\".\nReplace-UniqueText \"arteriovenous fistula (AVF). Using Fast Fourier transform (FFT), my \" \"arteriovenous fistula (AVF). My \"\n\n# 9) Body paragraph: quantify the monthly savings.\nReplace-UniqueText \" and thus reduced the monthly cost for patients using our products.\" \" and thus reduced around `$900 monthly cost for patients using our products.\"\n\n# 10) Closing paragraph: drop the repeated \"Data Scientist\" qualifier.\nReplace-UniqueText \"this Data Scientist role\" \"this role\"\nReplace-UniqueText \"about the Data Scientist position and\" \"about the position and\"\n"}
